$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card20")

for ($r = 3; $r -le 12; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    # Force the assigned value to be stored as text (matching the existing
    # inline-string "2" -> "20" cell type) instead of Excel's default
    # numeric auto-detection, then strip the temporary Text number-format
    # so the cell's style/formatting stays exactly as it was before.
    $cell.NumberFormat = "@"
    $cell.Value = "20"
    $cell.ClearFormats()
}
